$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two new rows (20 and 21) following the same pattern as the
# preceding rows in the "PageInstituciones" data pool sheet.
$ws.Range("A20").Value = "DEC_0101"
$ws.Range("B20").Value = "13712759-8"
$ws.Range("C20").Value = "Verity1.0"
$ws.Range("D20").Value = "verity"

$ws.Range("A21").Value = "DEC_0102"
$ws.Range("B21").Value = "13712759-8"
$ws.Range("C21").Value = "Verity1.0"
$ws.Range("D21").Value = "verity"

# Update the active selection to match the new editing position.
$ws.Range("B26").Select()
